$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + date range) ---
$ws.Range("A8").Value2 = "Volume 31   Number  49"
$ws.Range("C9").Value2 = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Group B: numeric cell -> "no data" string placeholder ---
# Row 14 holds stable placeholder cells we reuse as copy-source:
#   C14/D14/F14/G14 = shared string "0" (style 13), E14/H14 = shared string "***.*" (style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4104) | Out-Null  # xlPasteAll: C15 -> "0" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4104) | Out-Null  # xlPasteAll: F18 -> "0" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4104) | Out-Null  # xlPasteAll: C20 -> "0" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4104) | Out-Null  # xlPasteAll: D20 -> "0" (value + style 13)
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4104) | Out-Null  # xlPasteAll: E20 -> "***.*" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4104) | Out-Null  # xlPasteAll: C22 -> "0" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4104) | Out-Null  # xlPasteAll: C27 -> "0" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4104) | Out-Null  # xlPasteAll: D28 -> "0" (value + style 13)
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4104) | Out-Null  # xlPasteAll: E28 -> "***.*" (value + style 13)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4104) | Out-Null  # xlPasteAll: D31 -> "0" (value + style 13)
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4104) | Out-Null  # xlPasteAll: E31 -> "***.*" (value + style 13)

# --- Group C: "no data" string placeholder -> numeric cell (restore numeric style) ---
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restore numeric style
$ws.Range("D17").Value2 = 5
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restore numeric style
$ws.Range("E17").Value2 = -20
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restore numeric style
$ws.Range("D25").Value2 = 2
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restore numeric style
$ws.Range("E25").Value2 = 100

# --- Group A: plain numeric value updates ---
# Row 16
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 6
$ws.Range("G16").Value2 = 14
$ws.Range("H16").Value2 = -57.142857142857
$ws.Range("I16").Value2 = 119
$ws.Range("J16").Value2 = 141
$ws.Range("K16").Value2 = -15.602836879432
$ws.Range("L16").Value2 = 17.821782178217
$ws.Range("M16").Value2 = -24.683544303797
$ws.Range("N16").Value2 = -83.653846153846
# Row 17
$ws.Range("C17").Value2 = 4
$ws.Range("F17").Value2 = 13
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = -13.333333333333
$ws.Range("I17").Value2 = 170
$ws.Range("J17").Value2 = 174
$ws.Range("K17").Value2 = -2.298850574712
$ws.Range("L17").Value2 = 22.302158273381
$ws.Range("M17").Value2 = 65.048543689320
$ws.Range("N17").Value2 = -39.501779359430
# Row 18
$ws.Range("H18").Value2 = -100
$ws.Range("J18").Value2 = 90
$ws.Range("K18").Value2 = -47.777777777777
$ws.Range("L18").Value2 = -58.771929824561
$ws.Range("M18").Value2 = -38.157894736842
$ws.Range("N18").Value2 = -91.854419410745
# Row 19
$ws.Range("C19").Value2 = 7
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = -22.222222222222
$ws.Range("F19").Value2 = 30
$ws.Range("G19").Value2 = 32
$ws.Range("H19").Value2 = -6.25
$ws.Range("I19").Value2 = 331
$ws.Range("J19").Value2 = 372
$ws.Range("K19").Value2 = -11.021505376344
$ws.Range("L19").Value2 = 2.795031055900
$ws.Range("M19").Value2 = 30.314960629921
$ws.Range("N19").Value2 = -47.874015748031
# Row 20
$ws.Range("I20").Value2 = 42
$ws.Range("K20").Value2 = -58
$ws.Range("L20").Value2 = -45.454545454545
$ws.Range("M20").Value2 = 61.538461538461
$ws.Range("N20").Value2 = -87.719298245614
# Row 21
$ws.Range("C21").Value2 = 12
$ws.Range("D21").Value2 = 17
$ws.Range("E21").Value2 = -29.411764705882
$ws.Range("G21").Value2 = 72
$ws.Range("H21").Value2 = -25
$ws.Range("I21").Value2 = 719
$ws.Range("J21").Value2 = 886
$ws.Range("K21").Value2 = -18.848758465011
$ws.Range("L21").Value2 = -6.623376623376
$ws.Range("M21").Value2 = 13.765822784810
$ws.Range("N21").Value2 = -72.346153846153
# Row 22
$ws.Range("D22").Value2 = 1
$ws.Range("E22").Value2 = -100
$ws.Range("F22").Value2 = 3
$ws.Range("G22").Value2 = 4
$ws.Range("H22").Value2 = -25
$ws.Range("J22").Value2 = 36
$ws.Range("K22").Value2 = -27.777777777777
# Row 23
$ws.Range("D23").Value2 = 2
$ws.Range("E23").Value2 = -50
$ws.Range("F23").Value2 = 7
$ws.Range("G23").Value2 = 9
$ws.Range("H23").Value2 = -22.222222222222
$ws.Range("I23").Value2 = 123
$ws.Range("J23").Value2 = 143
$ws.Range("K23").Value2 = -13.986013986014
$ws.Range("L23").Value2 = 17.142857142857
$ws.Range("M23").Value2 = 64
# Row 24
$ws.Range("C24").Value2 = 10
$ws.Range("D24").Value2 = 9
$ws.Range("E24").Value2 = 11.111111111111
$ws.Range("G24").Value2 = 37
$ws.Range("H24").Value2 = 5.405405405405
$ws.Range("I24").Value2 = 396
$ws.Range("J24").Value2 = 479
$ws.Range("K24").Value2 = -17.327766179540
$ws.Range("L24").Value2 = -17.5
$ws.Range("M24").Value2 = -31.25
# Row 25
$ws.Range("C25").Value2 = 4
$ws.Range("F25").Value2 = 10
$ws.Range("G25").Value2 = 6
$ws.Range("H25").Value2 = 66.666666666666
$ws.Range("I25").Value2 = 75
$ws.Range("J25").Value2 = 149
$ws.Range("K25").Value2 = -49.664429530201
$ws.Range("L25").Value2 = -48.979591836734
# Row 26
$ws.Range("D26").Value2 = 4
$ws.Range("E26").Value2 = 25
$ws.Range("F26").Value2 = 25
$ws.Range("G26").Value2 = 20
$ws.Range("H26").Value2 = 25
$ws.Range("I26").Value2 = 235
$ws.Range("J26").Value2 = 228
$ws.Range("K26").Value2 = 3.070175438596
$ws.Range("L26").Value2 = 22.395833333333
$ws.Range("M26").Value2 = -15.467625899280
# Row 28
$ws.Range("C28").Value2 = 2
$ws.Range("F28").Value2 = 5
$ws.Range("G28").Value2 = 2
$ws.Range("H28").Value2 = 150
$ws.Range("I28").Value2 = 41
$ws.Range("K28").Value2 = -6.818181818181
$ws.Range("L28").Value2 = -21.153846153846
# Row 29
$ws.Range("N29").Value2 = -79.411764705882
# Row 30
$ws.Range("N30").Value2 = -79.310344827586
# Row 31
$ws.Range("G31").Value2 = 2
